$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E (price / volume) keep their original text formatting so
# values such as "1.00" or "6.00" are not coerced into plain numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '61.559.60'
$ws.Range('E2').Value = '  -2.21%  '
$ws.Range('D3').Value = '3.000.26'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '540.30'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').Value = '130.96'
$ws.Range('E6').Value = '  -5.02%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '2.995.18'
$ws.Range('E8').Value = '  -2.14%  '
$ws.Range('D9').Value = '0.488'
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('D10').Value = '6.00'
$ws.Range('E10').Value = '  -4.30%  '
$ws.Range('D11').Value = '0.144'
$ws.Range('E11').Value = '  -7.67%  '
$ws.Range('D12').Value = '0.444'
$ws.Range('E12').Value = '  -2.26%  '
$ws.Range('D13').Value = '34.13'
$ws.Range('E13').Value = '  -1.22%  '
$ws.Range('E14').Value = '  -2.22%  '
$ws.Range('D15').Value = '3.492.79'
$ws.Range('E15').Value = '  -2.23%  '
$ws.Range('D16').Value = '61.703.52'
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('D17').Value = '0.109'
$ws.Range('E17').Value = '  -2.95%  '
$ws.Range('D18').Value = '3.001.73'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('D19').Value = '6.59'
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').Value = '482.74'
$ws.Range('E20').Value = '  +2.35%  '
$ws.Range('D21').Value = '13.17'
$ws.Range('E21').Value = '  -2.41%  '
$ws.Range('D22').Value = '0.664'
$ws.Range('E22').Value = '  -4.83%  '
$ws.Range('D23').Value = '6.95'
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('D24').Value = '82.09'
$ws.Range('E24').Value = '  +4.60%  '
$ws.Range('D25').Value = '11.95'
$ws.Range('E25').Value = '  -1.58%  '
$ws.Range('D27').Value = '2.68'
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('D28').Value = '7.61'
$ws.Range('E28').Value = '  -3.40%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '1.90'
$ws.Range('E30').Value = '  +1.06%  '
$ws.Range('D31').Value = '25.55'
$ws.Range('E31').Value = '  -2.10%  '
$ws.Range('E32').Value = '  -3.07%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').Value = '5.55'
$ws.Range('E33').Value = '  +1.37%  '
$ws.Range('B34').Value = 'Stacks'
$ws.Range('C34').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D34').Value = '2.31'
$ws.Range('E34').Value = '  +0.99%  '
$ws.Range('D35').Value = '54.89'
$ws.Range('E35').Value = '  -6.64%  '
$ws.Range('D36').Value = '5.82'
$ws.Range('E36').Value = '  -2.79%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '3.134.49'
$ws.Range('E37').Value = '  -3.62%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '436.50'
$ws.Range('E38').Value = '  -9.73%  '
$ws.Range('D39').Value = '0.0791'
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('D40').Value = '0.0381'
$ws.Range('E40').Value = '  -4.50%  '
$ws.Range('D41').Value = '0.117'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('D42').Value = '8.04'
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('D43').Value = '2.38'
$ws.Range('E43').Value = '  -7.37%  '
$ws.Range('D45').Value = '26.06'
$ws.Range('E45').Value = '  +3.27%  '
$ws.Range('D46').Value = '0.241'
$ws.Range('E46').Value = '  -4.25%  '
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '115.82'
$ws.Range('E48').Value = '  -5.32%  '
$ws.Range('B49').Value = 'Fetch.AI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D49').Value = '1.93'
$ws.Range('E49').Value = '  -4.09%  '
$ws.Range('D50').Value = '1.27'
$ws.Range('E50').Value = '  +3.94%  '
$ws.Range('D51').Value = '0.0₃0479'
$ws.Range('E51').Value = '  -8.14%  '
